$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values (column B) per the diff
$ws.Range("B4").Value = 219.44          # Fuel/oxidizer mass (kg)
$ws.Range("B7").Value = 0               # Ullage volume (-)
$ws.Range("B8").Value = 1.05            # Proof factor (-)
$ws.Range("B10").Value = 16.9           # L/D (-)
$ws.Range("B13").Value = "18,66,86"     # Weave pattern (degrees)
$ws.Range("B15").Value = 0.28           # Ply thickness (mm)
$ws.Range("B25").Value = 260            # Liner yield strength (MPa)
$ws.Range("B26").Value = 72             # Liner E modulus (GPa)

# Update view/selection state to match the saved workbook
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("E9").Select()
